$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row for "Andrew Kaul" (special slot) right after Ann Gawel (row 9) ---
# This pushes the existing rows 10-18 down to 11-19, matching the target layout where
# the new entry becomes row 10.
$ws.Rows(10).Insert()

# New row 10: Andrew Kaul, "special" entry, same day as the Nov-29 session, 10:20 AM slot.
$ws.Range("A10").Value = "Andrew"
$ws.Range("B10").Value = "Kaul"
$ws.Range("C10").Value = "special"
$ws.Range("D10").Value = 42703
$ws.Range("E10").Value = 0.43055555555555558

# --- Re-sort the remainder: "Karri Folks" (now at row 14, random number 12) moves to the
# very end of the list (new last row, random number 17); "Andrea Rabinowitz" (was the
# last row, random number 17) takes over Karri's old slot (random number 12), and every
# row between keeps its day/time slot but is renumbered down by one. ---
$ws.Range("A14").Value = "Andrea              "
$ws.Range("B14").Value = "Rabinowitz                    "
$ws.Range("C14").Value = 12

$ws.Range("A15").Value = "Alexander           "
$ws.Range("B15").Value = "Karnish                       "
$ws.Range("C15").Value = 13

$ws.Range("A16").Value = "Amy                 "
$ws.Range("B16").Value = "Geffre                        "
$ws.Range("C16").Value = 14

$ws.Range("A17").Value = "Clare               "
$ws.Range("B17").Value = "Adams                         "
$ws.Range("C17").Value = 15

$ws.Range("A18").Value = "Hilary"
$ws.Range("B18").Value = "Haley"
$ws.Range("C18").Value = 16

$ws.Range("A19").Value = "Karri"
$ws.Range("B19").Value = "Folks"
$ws.Range("C19").Value = 17

# Blank trailing row (row 20), keeping the same time-column formatting as the row above.
$ws.Range("E20").NumberFormat = $ws.Range("E19").NumberFormat

$ws.Range("A20").Select()
